# Ficha workbook update: add a new "Especialidad" column, refresh the
# ficha numbers, and drop the thin cell borders (keeping the header fill).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Ficha numbers in column A -------------------------------
$ws.Range("A2").Value = 6565414
$ws.Range("A3").Value = 5441655
$ws.Range("A4").Value = 3325232

# --- Drop the thin outline border on the existing table, keep the fill -
$ws.Range("A1:B4").Borders.LineStyle = -4142

# --- Add the new "Especialidad" column (C) ------------------------------
# Copy the existing header / data formatting (now border-free) onto the
# new column so the fill colour + alignment match the rest of the table,
# then fill in the values (typed in the same order the original author
# entered them, so the shared-string table comes out in the same order).
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("C2:C4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C1").Value = "Especialidad"
$ws.Range("C3").Value = "lgtc"
$ws.Range("C2").Value = "adso"
$ws.Range("C4").Value = "adsi"

$ws.Columns.Item(3).ColumnWidth = 14.6

# --- Selection moves down to A6 after the edits -------------------------
$ws.Range("A6").Select()
